# Word COM-interop script implementing the diff:
#   1. Merge the two "SAT Dec 9" / " 10:13:13 PST 2017" runs into a single run.
#   2. After the *last* "Amount Received mode ... - CASH" paragraph, insert a
#      whole new purchase-record block (13 new paragraphs) before the blank
#      paragraph that already followed it.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge "SAT Dec 9" + " 10:13:13 PST 2017" into a single run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("SAT Dec 9 10:13:13 PST 2017", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "SAT Dec 9 10:13:13 PST 2017", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Locate the *last* occurrence of the "Amount Received mode ... - CASH"
#    paragraph (the one that sits right before the GAVI block end).
# ---------------------------------------------------------------------
$tab = [char]9
$searchText = "Amount Received mode" + $tab + $tab + "- CASH"

$scan = $d.Content
$scan.Start = 0
$lastStart = -1
$lastEnd = -1
while ($true) {
    $found = $scan.Find.Execute($searchText, $false, $false, $false, $false, $false, `
                                 $true, 1, $false, "", 0)
    if (-not $found) { break }
    $lastStart = $scan.Start
    $lastEnd = $scan.End
    $scan.Start = $scan.End
    $scan.End = $d.Content.End
}

# Translate the character offset into a document-level paragraph index
# (Range.Paragraphs.Item(1).Index is relative to the sub-range, not the
# document, so walk $d.Paragraphs instead).
$probe = $lastStart + 1
$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $probe -and $p.Range.End -gt $probe) {
        $anchorIndex = $i
        break
    }
}

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------
function New-BlankParagraphAfter($para) {
    $para.Range.InsertParagraphAfter()
    return $d.Paragraphs.Item($para.Index + 1)
}

function Set-ParaText($para, [string[]]$pieces, [bool]$bold) {
    # $pieces is an ordered list of tokens: plain strings are literal text,
    # the literal tab character denotes a <w:tab/> element.
    $rng = $para.Range
    if ($bold) { $rng.Bold = 1 }

    $insertAt = $para.Range.End - 1   # collapsed point just before the paragraph mark

    foreach ($piece in $pieces) {
        $ip = $d.Range($insertAt, $insertAt)
        if ($piece -eq $tab) {
            $ip.InsertBefore($tab)
        } else {
            $ip.InsertBefore($piece)
        }
        $insertAt = $para.Range.End - 1
    }
}

$anchorPara = $d.Paragraphs.Item($anchorIndex)

# --- Paragraph 1: blank line -------------------------------------------------
$p1 = New-BlankParagraphAfter $anchorPara

# --- Paragraph 2: "SUN Dec 10" + " 10:05:04 PST 2017" (kept as two runs) ----
$p2 = New-BlankParagraphAfter $p1
$ip1 = $d.Range($p2.Range.End - 1, $p2.Range.End - 1)
$ip1.InsertBefore("SUN Dec 10")
$ip2 = $d.Range($p2.Range.End - 1, $p2.Range.End - 1)
$ip2.InsertBefore(" 10:05:04 PST 2017")
# Force the two pieces to stay as separate runs (same visible formatting,
# but toggling Bold on just the 2nd piece stops the engine auto-merging it
# back into the first run) while leaving no stray formatting override behind.
$splitStart = $p2.Range.Start + ("SUN Dec 10").Length
$splitEnd = $p2.Range.End - 1
$splitRng = $d.Range($splitStart, $splitEnd)
$splitRng.Bold = 1
$splitRng.Bold = 0

# --- Paragraph 3: "Person Name" ... "- GAVI" --------------------------------
$p3 = New-BlankParagraphAfter $p2
Set-ParaText $p3 @("Person Name", $tab, $tab, $tab, $tab, "- GAVI") $false

# --- Paragraph 4: "Bill number" ... "- 1866" --------------------------------
$p4 = New-BlankParagraphAfter $p3
Set-ParaText $p4 @("Bill number", $tab, $tab, $tab, $tab, "- 1866") $false

# --- Paragraph 5: dashed separator ------------------------------------------
$p5 = New-BlankParagraphAfter $p4
Set-ParaText $p5 @("---------------------------------------------------------------") $false

# --- Paragraph 6: "Item Name" ... "- PUDI CARROT" ---------------------------
$p6 = New-BlankParagraphAfter $p5
Set-ParaText $p6 @("Item Name", $tab, $tab, $tab, $tab, "- PUDI CARROT") $false

# --- Paragraph 7: "Number of Pockets" ... "- 1" -----------------------------
$p7 = New-BlankParagraphAfter $p6
Set-ParaText $p7 @("Number of Pockets", $tab, $tab, $tab, "- 1") $false

# --- Paragraph 8: "Number of KGs" ... "- 100" -------------------------------
$p8 = New-BlankParagraphAfter $p7
Set-ParaText $p8 @("Number of KGs", $tab, $tab, $tab, "- 100") $false

# --- Paragraph 9: "Rate" ... "- 12" -----------------------------------------
$p9 = New-BlankParagraphAfter $p8
Set-ParaText $p9 @("Rate", $tab, $tab, $tab, $tab, $tab, "- 12") $false

# --- Paragraph 10: "Total Price" ... "- 1200.0" -----------------------------
$p10 = New-BlankParagraphAfter $p9
Set-ParaText $p10 @("Total Price", $tab, $tab, $tab, $tab, "- 1200.0") $false

# --- Paragraph 11: "Amount balance" (bold) ... "- 20014.5" ------------------
$p11 = New-BlankParagraphAfter $p10
Set-ParaText $p11 @("Amount balance", $tab, $tab, $tab, "- 20014.5") $true

# --- Paragraph 12: blank line -----------------------------------------------
$p12 = New-BlankParagraphAfter $p11
